# chore: rename package to `player`
#
# Applies the logo.pptx edit:
#   1. Resize/move the dark "Rectangle 9" backing panel behind the controller art.
#   2. Remove the extra, now-unused "Picture 4" (duplicate faded controller image).
#   3. Nudge the "Picture 7" controller picture up/left slightly.
#   4. Add a new "PLAYER" title text box (SNES-styled) over the artwork.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Rectangle 9 ("Rectangle 8"-style backing panel) — move + grow taller.
#    (Point values below are chosen so the saved EMUs land exactly on
#    8832527,9471467 / 3887078,5797866 — the rotation (20627742) is untouched.)
$rect9 = $s.Shapes.Item("Rectangle 9")
$rect9.Left = 695.4745869291338
$rect9.Top = 745.7848231496063
$rect9.Width = 306.0691438582677
$rect9.Height = 456.52489188976375

# 2) Drop the redundant faded controller picture entirely.
$s.Shapes.Item("Picture 4").Delete()

# 3) Picture 7 (controller render) — reposition only; size/rotation unchanged.
$pic7 = $s.Shapes.Item("Picture 7")
$pic7.Left = 723.6656892913385
$pic7.Top = 1004.863937007874

# 4) New "PLAYER" text box, rotated with the art, SNES-styled.
$textBox = $s.Shapes.AddTextbox(1, 100, 100, 200, 50)
$tr = $textBox.TextFrame.TextRange
$tr.Text = "PLAYER"
$tr.ParagraphFormat.Alignment = 2
$tr.Font.Size = 220
$tr.Font.Name = "SNES"
$tr.Font.Color.RGB = 5132874

$textBox.TextFrame.WordWrap = $true
$textBox.TextFrame.AutoSize = 1
$textBox.Fill.Visible = $false

$textBox.Left = 556.8232283464567
$textBox.Top = 680.8545869291338
$textBox.Width = 492.57771653543307
$textBox.Height = 273.84843519685035
$textBox.Rotation = 343.2283
